$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the two new use-case rows in column E first, so the shared-string
# table interns them ahead of the REQ(13)/REQ(14) text (matches author order).
$ws.Range("E23").Value = "UC10 - Recolher Numerário"
$ws.Range("E24").Value = "UC11 - Emitir Ingresso"
$ws.Rows.Item(23).RowHeight = 12.75
$ws.Rows.Item(24).RowHeight = 12.75

# Insert the two new requirement rows (13 and 14) with their three columns.
$ws.Range("A13").Value = "REQ(13) "
$ws.Range("B13").Value = "Após todo o processo feito com sucesso, o cliente recebeo ticket para assistir o filme (Cliente)"
$ws.Range("C13").Value = "Emitir Ingresso"

$ws.Range("A14").Value = "REQ(14)"
$ws.Range("B14").Value = "Após várias vendas acumuladas, o responsável pela máquina e pelo sistema, vai até as máquinas retirar as quantias (Gerente)"
$ws.Range("C14").Value = "Recolher Numerário"

# Widen column B to fit the new (longer) content.
# (Target stored width is 108.85546875; this runtime's ColumnWidth setter
# quantises to the nearest 1/6 + 5/6 pixel step, so 108 is the closest input
# that reproduces the saved width, 108.83333333333333.)
$ws.Columns.Item(2).ColumnWidth = 108

# Match the saved selection state.
$ws.Range("C15").Select()
